$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.271826505661011
$ws.Range("B1").Value = 4.877838611602783
$ws.Range("C1").Value = 3.216334581375122
$ws.Range("D1").Value = 2.448996067047119
$ws.Range("E1").Value = 2.015807867050171
